# Weekly refresh: insert 2 new observation rows (240, 241) into the
# "Feria Lagunitas de Puerto Montt - Cilantro" daily logic sheet.
# Every existing data row at/after 240 shifts down by 2 (handled by
# Rows.Insert, which also carries the row-240 formatting, incl. the date
# number format on column D, onto the freshly inserted rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 240 (pushes old 240.. down to 242..)
$ws.Rows.Item(240).Insert()
$ws.Rows.Item(240).Insert()

# Static columns shared by every row in this sheet (single market / product)
$colA = 4
$colB = 'Feria Lagunitas de Puerto Montt'
$colC = 'Los Lagos'
$colE = 10
$colF = 100112040
$colG = 'Cilantro'
$colH = 'Sin especificar'
$colI = 'Primera'
$colR = 'Hortaliza'

function Set-DataRow($RowIndex, $D, $J, $K, $L, $M, $N, $O, $P, $Q) {
    $ws.Cells.Item($RowIndex, 1).Value = $colA
    $ws.Cells.Item($RowIndex, 2).Value = $colB
    $ws.Cells.Item($RowIndex, 3).Value = $colC
    $ws.Cells.Item($RowIndex, 4).Value = $D
    $ws.Cells.Item($RowIndex, 5).Value = $colE
    $ws.Cells.Item($RowIndex, 6).Value = $colF
    $ws.Cells.Item($RowIndex, 7).Value = $colG
    $ws.Cells.Item($RowIndex, 8).Value = $colH
    $ws.Cells.Item($RowIndex, 9).Value = $colI
    $ws.Cells.Item($RowIndex, 10).Value = $J
    $ws.Cells.Item($RowIndex, 11).Value = $K
    $ws.Cells.Item($RowIndex, 12).Value = $L
    $ws.Cells.Item($RowIndex, 13).Value = $M
    $ws.Cells.Item($RowIndex, 14).Value = $N
    $ws.Cells.Item($RowIndex, 15).Value = $O
    $ws.Cells.Item($RowIndex, 16).Value = $P
    $ws.Cells.Item($RowIndex, 17).Value = $Q
    $ws.Cells.Item($RowIndex, 18).Value = $colR
}

# New row 240: 2022-06-10, $/caja 36 atados, Región Metropolitana
Set-DataRow 240 44722 150 11000 11000 11000 '$/caja 36 atados' 'Región Metropolitana' 306 36

# New row 241: 2022-06-10, $/docena de atados (2 kilos), Región de La Araucanía
Set-DataRow 241 44722 120 5000 5000 5000 '$/docena de atados (2 kilos)' 'Región de La Araucanía' 2500 2
